# Apply the "add rusentiment / rutweetcorp evaluation rows" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Template rows for styling -------------------------------------------------
# Row 19 carries the "plain" style group (s=0 for text/int cells, s=2 for decimals).
# Row 15 carries the "fill" style group (s=5 for text/int cells, s=6 for decimals).
$plainTemplateRow = 19
$fillTemplateRow   = 15

# --- New data rows (dataset, f_macro_2, f_macro_3, f_micro_2, f_micro_3, f_neg, f_neutral, f_pos, negative, neutral, positive, model_name) ---
$newRows = @(
    @{ Row = 20; Template = $plainTemplateRow; A = "rusentiment"; B = 0.533225755251045;  C = 0.592862047640053;  D = 0.537313432835821;  E = 0.636;  F = 0.463414634146341;  G = 0.712134632418069;  H = 0.603036876355748;  I = 189; J = 587; K = 224; L = "dostoevsky" },
    @{ Row = 21; Template = $plainTemplateRow; A = "rusentiment"; B = 0.571727290410426;  C = 0.665913431702189;  D = 0.59;               E = 0.775;  F = 0.509433962264151;  G = 0.854285714285714;  H = 0.634020618556701;  I = 100; J = 726; K = 174; L = "dostoevsky" },
    @{ Row = 22; Template = $plainTemplateRow; A = "rusentiment"; B = 0.337668643121079;  C = 0.225112428747386;  D = 0.339514978601997;  E = 0.238;  F = 0.345707656612529;  G = 0;                   H = 0.32962962962963;   I = 184; J = 598; K = 218; L = "sentiment_twitter_preproc" },
    @{ Row = 23; Template = $plainTemplateRow; A = "rusentiment"; B = 0.254342021614749;  C = 0.169561347743166;  D = 0.253386454183267;  E = 0.159;  F = 0.227692307692308;  G = 0;                   H = 0.28099173553719;   I = 85;  J = 745; K = 170; L = "sentiment_twitter_preproc" },
    @{ Row = 24; Template = $fillTemplateRow;  A = "rusentiment"; B = 0.644455212515575;  C = 0.679714762818304;  D = 0.644468313641246;  E = 0.701;  F = 0.632258064516129;  G = 0.75023386342376;   H = 0.656652360515021;  I = 195; J = 580; K = 225; L = "rusentiment_convers_bert" },
    @{ Row = 25; Template = $fillTemplateRow;  A = "rusentiment"; B = 0.702954561872426;  C = 0.764615437155969;  D = 0.704507512520868;  E = 0.833;  F = 0.697247706422018;  G = 0.887937187723055;  H = 0.708661417322835;  I = 101; J = 720; K = 179; L = "rusentiment_convers_bert" },
    @{ Row = 26; Template = $fillTemplateRow;  A = "rutweetcorp"; B = 0.708693599681972;  C = 0.699160662031267;  D = 0.737024221453287;  E = 0.713;  F = 0.857558139534884;  G = 0.680094786729858;  H = 0.55982905982906;   I = 330; J = 352; K = 318; L = "rusentiment_convers_bert" },
    @{ Row = 27; Template = $plainTemplateRow; A = "rutweetcorp"; B = 0.683858441392601;  C = 0.455905627595067;  D = 0.678271308523409;  E = 0.565;  F = 0.75130890052356;   G = 0;                   H = 0.616407982261641;  I = 330; J = 334; K = 336; L = "sentiment_twitter_preproc" },
    @{ Row = 28; Template = $plainTemplateRow; A = "rutweetcorp"; B = 0.512053707659445;  C = 0.53936913843963;   D = 0.532;               E = 0.563;  F = 0.665486725663717;  G = 0.594;               H = 0.358620689655172;  I = 354; J = 308; K = 338; L = "dostoevsky" }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")

foreach ($rowSpec in $newRows) {
    $targetRow = $rowSpec.Row
    $templateRow = $rowSpec.Template

    # Copy the template row (formatting + number formats) down onto the new row first.
    $srcRange = $ws.Range("A" + $templateRow + ":L" + $templateRow)
    $srcRange.Copy()
    $dstRange = $ws.Range("A" + $targetRow + ":L" + $targetRow)
    $dstRange.PasteSpecial(-4122) # xlPasteFormats

    foreach ($col in $cols) {
        $cellAddr = $col + $targetRow
        $ws.Range($cellAddr).Value = $rowSpec[$col]
    }
}

# --- Misc sheet metadata tweaks from the diff ----------------------------------
# (Column-width / defaultColWidth deltas in the diff are sub-pixel cosmetic
#  artifacts from the original authoring tool's relayout; the COM width
#  setters here only resolve to whole-pixel buckets, so touching them would
#  move the stored value further from the target than leaving it alone.)
$ws.Range("B9").Select()
